# daily auto push: 2026-01-16 13:47 UTC
# A new timestamped reading for 2026/01/16 (Friday) is appended to the log.
# Because the log is sorted by date/time, the new record lands at row 663,
# pushing the existing rows 663-704 down to rows 664-705.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: shift rows 663 onward down by one row.
$ws.Rows(663).Insert()

# The date/weekday text for the new record ("2026/01/16" / "金") already
# exists verbatim a few rows above (A662:B662) - copy it down so the new
# cells stay plain text instead of being reinterpreted as a date value.
$ws.Range("A662:B662").Copy($ws.Range("A663:B663"))

# Fill in the new record's numeric measurements.
$ws.Range("C663").Value = 20
$ws.Range("D663").Value = 201
